$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "884D6D"
$ws.Range("B6").Value = "Almohadilla + Chip Epson T04D1"
$ws.Range("C6").Value = "ET M1140 M1180 M2140 M3140 M3170 M3180 2700 2750 2756 2760 3700 3710 3750 3760 L4100 L4150 L4151 L4153 L4156 L4158 L4160 L4163 L4165 L4166 L4167 L4168 L4169 L4170 L4260 4700 4750 4760 L6160 L6161 L6168 L6170 L6171 L6178 L6190 L6191 L6198 L14150, WF 2800 2860 2865 ST3000 ST4000, XP 5100 5105 5115 8500"
$ws.Range("D6").Value = 35000
$ws.Range("E6").Value = 200000
$ws.Range("F6").Value = 10
$ws.Range("G6").Value = 18
$ws.Range("H6").Formula = "=(E6-D6)*G6"
$ws.Range("I6").Formula = "=D6*F6"
$ws.Range("J6").Value = 350000
